# Update cryptos list: refresh Price (D) and Volume(1h) (E) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.940.35"
$ws.Range("E2").Value = "  -0.59%  "
$ws.Range("D3").Value = "2.094.71"
$ws.Range("E3").Value = "  +2.00%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "'245.84"
$ws.Range("E5").Value = "  -0.99%  "
$ws.Range("D6").Value = "'0.655"
$ws.Range("E6").Value = "  -1.57%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "'56.07"
$ws.Range("E8").Value = "  -3.71%  "
$ws.Range("D9").Value = "'60.02"
$ws.Range("E9").Value = "  -0.14%  "
$ws.Range("D10").Value = "'0.370"
$ws.Range("E10").Value = "  -3.29%  "
$ws.Range("E11").Value = "  -0.76%  "
$ws.Range("E12").Value = "  +1.31%  "
$ws.Range("D13").Value = "'15.19"
$ws.Range("E13").Value = "  -4.67%  "
$ws.Range("E14").Value = "  +7.06%  "
$ws.Range("D15").Value = "2.405.29"
$ws.Range("E15").Value = "  +2.17%  "
$ws.Range("E16").Value = "  -2.61%  "
$ws.Range("D17").Value = "2.152.02"
$ws.Range("E17").Value = "  +4.72%  "
$ws.Range("D18").Value = "36.883.61"
$ws.Range("E18").Value = "  -0.87%  "
$ws.Range("D19").Value = "'17.49"
$ws.Range("E19").Value = "  -3.69%  "
$ws.Range("E20").Value = "  -1.57%  "
$ws.Range("E21").Value = "  -0.99%  "
$ws.Range("D22").Value = "'5.52"
$ws.Range("E22").Value = "  +3.07%  "
$ws.Range("D23").Value = "'238.44"
$ws.Range("E23").Value = "  +0.49%  "
$ws.Range("D24").Value = "'1.00"
$ws.Range("E24").Value = "  -0.02%  "
$ws.Range("D25").Value = "'2.43"
$ws.Range("E25").Value = "  -1.47%  "
$ws.Range("D26").Value = "'9.91"
$ws.Range("E26").Value = "  +5.08%  "
$ws.Range("E27").Value = "  -0.09%  "
$ws.Range("D28").Value = "'168.81"
$ws.Range("E28").Value = "  -0.40%  "
$ws.Range("D29").Value = "'20.84"
$ws.Range("E29").Value = "  +3.93%  "
$ws.Range("D30").Value = "'5.42"
$ws.Range("E30").Value = "  +13.15%  "
$ws.Range("E31").Value = "  -0.28%  "
$ws.Range("D32").Value = "'1.20"
$ws.Range("E32").Value = "  +6.65%  "
$ws.Range("D33").Value = "'4.71"
$ws.Range("E33").Value = "  +4.95%  "
$ws.Range("D34").Value = "'0.0613"
$ws.Range("E34").Value = "  -0.89%  "
$ws.Range("D35").Value = "'2.41"
$ws.Range("E35").Value = "  +6.59%  "
$ws.Range("E36").Value = "  +0.10%  "
$ws.Range("E37").Value = "  +4.26%  "
$ws.Range("D38").Value = "'0.0843"
$ws.Range("E38").Value = "  -6.66%  "
$ws.Range("E39").Value = "  -3.69%  "
$ws.Range("E40").Value = "  +2.29%  "
$ws.Range("E41").Value = "  -0.14%  "
$ws.Range("E42").Value = "  -5.89%  "
$ws.Range("D43").Value = "'0.0957"
$ws.Range("E43").Value = "  -6.40%  "
$ws.Range("E44").Value = "  +1.51%  "
$ws.Range("D45").Value = "'2.87"
$ws.Range("E45").Value = "  -11.28%  "
$ws.Range("D46").Value = "'16.17"
$ws.Range("E46").Value = "  -6.08%  "
$ws.Range("D47").Value = "1.357.43"
$ws.Range("E47").Value = "  +6.40%  "
$ws.Range("D48").Value = "'2.48"
$ws.Range("E48").Value = "  +0.75%  "
$ws.Range("E49").Value = "  +3.26%  "
$ws.Range("E50").Value = "  +0.09%  "
$ws.Range("D51").Value = "2.287.76"
$ws.Range("E51").Value = "  +2.10%  "
